# Update vm_pu results for Case_5_49 (380 kV slack voltage set-point case)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row number -> ordered column letters -> new numeric value
$updates = @{
    2 = @{ "B"=1.02; "C"=1.039690384347496; "D"=1.043064676187131; "E"=1.047083747646934; "F"=1.055334453513962; "I"=1.041427754092698; "J"=1.044781034722957; "K"=1.045839595071604; "L"=1.04984737626476; "M"=1.058075196496949 }
    3 = @{ "B"=1.02; "C"=1.040609577373308; "D"=1.043770930488108; "E"=1.047935789649966; "F"=1.056349101702797; "I"=1.041679943851032; "J"=1.045345689014041; "K"=1.046357082479689; "L"=1.050511086525433; "M"=1.058902752084843 }
    4 = @{ "B"=1.02; "C"=1.041204687713334; "D"=1.044228070383623; "E"=1.048487844858186; "F"=1.057006714139265; "I"=1.041841896329369; "J"=1.04571075354616; "K"=1.046691403682436; "L"=1.050940633937386; "M"=1.059438690554835 }
    5 = @{ "B"=1.02; "C"=1.041454950318724; "D"=1.044420285625038; "E"=1.04872010162592; "F"=1.05728342864312; "I"=1.041909685801266; "J"=1.045864152826447; "K"=1.046831824988231; "L"=1.051121234383256; "M"=1.059664106800137 }
    6 = @{ "B"=1.02; "C"=1.041496975068551; "D"=1.04445256134245; "E"=1.048759108652024; "F"=1.057329905097941; "I"=1.041921050614767; "J"=1.045889904882782; "K"=1.046855394848232; "L"=1.051151559060127; "M"=1.059701961462943 }
    7 = @{ "B"=1.02; "C"=1.041208031428442; "D"=1.044230638642651; "E"=1.048490947607293; "F"=1.057010410614608; "I"=1.041842803296893; "J"=1.045712803566203; "K"=1.04669328050034; "L"=1.050943047057147; "M"=1.059441702155922 }
    8 = @{ "B"=1.02; "C"=1.040000960786521; "D"=1.043303327061337; "E"=1.047371547906699; "F"=1.055677136890552; "I"=1.041513237355956; "J"=1.044971924779211; "K"=1.046014591048244; "L"=1.050071662551822; "M"=1.058354777957168 }
    9 = @{ "B"=1.02; "C"=1.037876528091758; "D"=1.041670464601232; "E"=1.045404648419296; "F"=1.05333596729872; "I"=1.040923097686924; "J"=1.043664110595814; "K"=1.04481465487762; "L"=1.048536849176095; "M"=1.05644301099452 }
    10 = @{ "B"=1.02; "C"=1.036462033688124; "D"=1.040582762277525; "E"=1.044097241174737; "F"=1.051780794084998; "I"=1.040523385298083; "J"=1.042790749141808; "K"=1.044012066693779; "L"=1.047514156549238; "M"=1.055170941308986 }
    11 = @{ "B"=1.02; "C"=1.035849979625534; "D"=1.040111998788371; "E"=1.043532050627647; "F"=1.051108731918158; "I"=1.040348821826855; "J"=1.042412233057824; "K"=1.043663925119564; "L"=1.047071454383241; "M"=1.054620713970781 }
    12 = @{ "B"=1.02; "C"=1.035622701221376; "D"=1.039937170454716; "E"=1.043322253981538; "F"=1.050859300329804; "I"=1.040283758576276; "J"=1.04227158440704; "K"=1.043534518279036; "L"=1.046907035672839; "M"=1.05441642433547 }
    13 = @{ "B"=1.02; "C"=1.035671450221496; "D"=1.039974670157616; "E"=1.043367249705812; "F"=1.050912795082026; "I"=1.040297724928803; "J"=1.042301756310342; "K"=1.043562280613952; "L"=1.046942303097064; "M"=1.054460241097607 }
    14 = @{ "B"=1.02; "C"=1.035831191375605; "D"=1.040097546718767; "E"=1.043514705902337; "F"=1.051088109693053; "I"=1.040343448216122; "J"=1.04240060803021; "K"=1.043653230172166; "L"=1.047057863055737; "M"=1.054603825497957 }
    15 = @{ "B"=1.02; "C"=1.035929621909786; "D"=1.04017325959296; "E"=1.043605577158775; "F"=1.051196153643307; "I"=1.040371590361817; "J"=1.042461507115794; "K"=1.043709255092001; "L"=1.047129066117352; "M"=1.054692304450726 }
    16 = @{ "B"=1.02; "C"=1.036502662867062; "D"=1.040614010021955; "E"=1.044134770614501; "F"=1.051825424931309; "I"=1.040534939253072; "J"=1.042815862818024; "K"=1.044035158812755; "L"=1.047543540066951; "M"=1.055207470549033 }
    17 = @{ "B"=1.02; "C"=1.036862232453386; "D"=1.040890540628318; "E"=1.04446696833425; "F"=1.052220509271005; "I"=1.040637006412749; "J"=1.043038049093456; "K"=1.044239425448929; "L"=1.047803564253176; "M"=1.055530778606474 }
    18 = @{ "B"=1.02; "C"=1.037072004866549; "D"=1.041051857274865; "E"=1.044660822849026; "F"=1.052451084309086; "I"=1.040696397106258; "J"=1.043167613260471; "K"=1.044358511263767; "L"=1.047955244486321; "M"=1.055719415313594 }
    19 = @{ "B"=1.02; "C"=1.037143538859232; "D"=1.041106865608258; "E"=1.044726937321498; "F"=1.052529726287538; "I"=1.040716623465784; "J"=1.043211785613355; "K"=1.044399106335839; "L"=1.048006965620622; "M"=1.055783745121485 }
    20 = @{ "B"=1.02; "C"=1.036823649734589; "D"=1.040860869327349; "E"=1.044431317426812; "F"=1.052178107105898; "I"=1.040626070396212; "J"=1.043014214049779; "K"=1.044217515715342; "L"=1.047775664816256; "M"=1.05549608485082 }
    21 = @{ "B"=1.02; "C"=1.035784149765376; "D"=1.04006136167984; "E"=1.043471279844989; "F"=1.051036478317765; "I"=1.04032998998815; "J"=1.042371500061194; "K"=1.043626450314154; "L"=1.047023832952025; "M"=1.054561540997144 }
    22 = @{ "B"=1.02; "C"=1.03513095560871; "D"=1.039558877934955; "E"=1.042868478399375; "F"=1.050319861584365; "I"=1.040142544836716; "J"=1.041967106074668; "K"=1.043254294765117; "L"=1.046551246005844; "M"=1.053974473361648 }
    23 = @{ "B"=1.02; "C"=1.035477189792414; "D"=1.039825234826177; "E"=1.043187957436618; "F"=1.05069964231641; "I"=1.040242034926423; "J"=1.042181510570191; "K"=1.043451631322228; "L"=1.046801761624401; "M"=1.054285639641836 }
    24 = @{ "B"=1.02; "C"=1.036841083476087; "D"=1.040874276446247; "E"=1.044447426261797; "F"=1.052197266420554; "I"=1.040631012354329; "J"=1.043024984183099; "K"=1.044227415964651; "L"=1.047788271332267; "M"=1.055511761290866 }
    25 = @{ "B"=1.02; "C"=1.038425433012668; "D"=1.042092450723107; "E"=1.045912464504228; "F"=1.053940232979099; "I"=1.041076773494662; "J"=1.044002477745126; "K"=1.045125335051925; "L"=1.048933548678692; "M"=1.056936823114824 }
}

foreach ($rowNum in $updates.Keys) {
    $rowData = $updates[$rowNum]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$rowNum").Value = $rowData[$col]
    }
}
